$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles/borders/fonts) from row 35, which has the same
# s="6"/s="7" style pattern that the new row 55 needs.
$ws.Range("A35:E35").Copy()
$ws.Range("A55:E55").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Match the row height used by row 35 and other similarly wrapped rows.
$ws.Rows.Item(55).RowHeight = 43.2

# Fill in the new row's values.
$ws.Range("A55").Value = "SCRIPT/T01P01A/us2301.ssb"
$ws.Range("B55").Value = 19
$ws.Range("C55").Value = " If you want to evolve, maybe\nyou should head over to [CS:P]Luminous\nSpring[CR] again."
$ws.Range("D55").Value = " Если вы захотите\nэволюционировать, вернитесь к [CS:P]Сияющему\nИсточнику[CR]."
$ws.Range("E55").Value = " Åòìé âú èàöïóéóå\nüâïìýøéïîéñïâàóû, âåñîéóåòû ë [CS:P]Òéÿýþåíô\nÉòóïœîéëô[CR]."

# Update the view so the new last row is visible and selected, mirroring
# the author's saved view state (scrolled down, D55 selected).
$excel.Goto($ws.Range("A53"), $true)
$ws.Range("D55").Select()
